$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that currently sits after the
#    "Coleten" run under "Super Jump" (it will be re-created below,
#    further down the document, next to the new "Level 9" text).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Split the "Level 9: Coleten" paragraph into three runs and add
#    the water blurb, re-inserting the _GoBack bookmark just before
#    the final " Coleten" run.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Level 9: Coleten`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Level 9: Coleten' paragraph"
}

$pStart = $target.Range.Start
# Range covering just the paragraph's own characters, not its
# trailing paragraph mark, so the paragraph's own properties /
# attributes are preserved by the replace.
$inner = $d.Range($pStart, $pStart + 16)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Level 9:</w:t></w:r><w:r><w:t xml:space="preserve"> Involves water somehow</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> Coleten</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$inner.InsertXML($xml)
